$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (B2/C2/D2) to their new contents
$ws.Range("B2").Value = "primer"
$ws.Range("C2").Value = "kickoff"
$ws.Range("D2").Value = "37'"

# Remove row 3 entirely, shifting the rest of the sheet up
$ws.Rows(3).Delete()
